$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("table_description")

# New variable names for the "reflectance" table block
$ws.Range("B57").Value = "ndvi"
$ws.Range("B58").Value = "notes"
$ws.Range("B59").Value = "pre_post_cut"

# Descriptions / units for the new rows
$ws.Range("C56").Value = "Time of sampling"
$ws.Range("D56").Value = "hh:mm"
$ws.Range("C57").Value = "NDVI value"
$ws.Range("C58").Value = "Notes"
$ws.Range("C59").Value = "Measurment was taken before or after the cut"

# TableID labels marking the start of each table block
$ws.Range("A32").Value = "cflux"
$ws.Range("A56").Value = "reflectance"

# Remaining cells that reuse already-existing shared strings
$ws.Range("B56").Value = "time"
$ws.Range("D57").Value = "percentage"
$ws.Range("E56").Value = "defined"
$ws.Range("E57").Value = "measured"
$ws.Range("D59").Value = "pre or post"
$ws.Range("E59").Value = "recorded"
